$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E6 value
$ws.Range("E6").Value = "28EhG4ya"

# Add rows 8, 9, 10 with "a" in columns A-E
foreach ($r in 8..10) {
    foreach ($c in 1..5) {
        $ws.Cells.Item($r, $c).Value = "a"
    }
}
